$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title text (shared string used by A1) - change year 2023 -> 2025
$ws.Range("A1").Value = "Total Demand for All Users customers for 2025 (Mthembanji)"

# Update monthly demand values (January..April rows)
$ws.Range("B3").Value = 695.3289142717998
$ws.Range("B4").Value = 562.3680155543
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
